# update v2 [new interface]
# Refreshes the sizing data in BASELINE_SIGLA and the rollup in Controle.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: BASELINE_SIGLA  (numeric STORAGE/MEMORIA/CPU columns stay numbers)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BASELINE_SIGLA")

$baseline = @(
    @("MA0283", "Homologation", "BACKUP 16GB RAM 16vCPU (WEBSERVER)", 230, 16, 16, "WEBSERVER", "Windows Server"),
    @("MA0284", "Homologation", "BACKUP 16GB RAM 16vCPU (WEBSERVER)", 230, 16, 16, "WEBSERVER", "Windows Server"),
    @("MA0285", "Homologation", "BACKUP 16GB RAM 16vCPU (WEBSERVER)", 170, 16, 16, "WEBSERVER", "Windows Server"),
    @("MA0286", "Homologation", "BACKUP 8GB RAM 4vCPU (WEBSERVER)",   170, 8,  4,  "WEBSERVER", "Windows Server"),
    @("MA0287", "Homologation", "BACKUP 16GB RAM 4vCPU (BACKUP)",     230, 16, 4,  "BACKUP",    "Windows Server"),
    @("MA0288", "Development",  "BACKUP 16GB RAM 8vCPU (BACKUP)",     170, 16, 8,  "BACKUP",    "Windows Server"),
    @("MA0289", "Development",  "BACKUP 8GB RAM 8vCPU (BACKUP)",      230, 8,  8,  "BACKUP",    "Windows Server")
)

for ($i = 0; $i -lt $baseline.Length; $i++) {
    $row = $i + 2
    $data = $baseline[$i]
    $ws1.Cells.Item($row, 1).Value = $data[0]
    $ws1.Cells.Item($row, 2).Value = $data[1]
    $ws1.Cells.Item($row, 3).Value = $data[2]
    $ws1.Cells.Item($row, 4).Value = $data[3]
    $ws1.Cells.Item($row, 5).Value = $data[4]
    $ws1.Cells.Item($row, 6).Value = $data[5]
    $ws1.Cells.Item($row, 7).Value = $data[6]
    $ws1.Cells.Item($row, 8).Value = $data[7]
}

# Drop the now-obsolete trailing rows (previously MA0297-MA0303)
$ws1.Range("A9:H15").Delete()

# ---------------------------------------------------------------------
# Sheet 2: Controle  (several "numeric-looking" columns are stored as TEXT
# in the source file, e.g. G="4" not 4. Excel would normally auto-coerce a
# numeric literal typed into a General cell into a real number, so we
# briefly mark the cell as Text before writing, then restore the original
# General/Normal styling so the stored style index is unaffected.)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Controle")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$controle = @(
    @("Development",  "4 x BACKUP 16GB RAM 8vCPU (BACKUP)",     "sa-east-1", "Windows Server", "c6i.2xlarge",  "Shared Instances", "4", "40", "Hours/Week", "On-Demand", "General Purpose SSD (gp3)", "230", "2x Daily", "10"),
    @("Homologation",  "6 x BACKUP 32GB RAM 16vCPU (WEBSERVER)", "sa-east-1", "Windows Server", "c6a.4xlarge",  "Shared Instances", "6", "40", "Hours/Week", "On-Demand", "General Purpose SSD (gp3)", "230", "2x Daily", "10"),
    @("Homologation",  "2 x BACKUP 8GB RAM 4vCPU (WEBSERVER)",   "sa-east-1", "Windows Server", "c6in.xlarge",  "Shared Instances", "2", "40", "Hours/Week", "On-Demand", "General Purpose SSD (gp3)", "170", "2x Daily", "10"),
    @("Homologation",  "2 x BACKUP 16GB RAM 4vCPU (BACKUP)",     "sa-east-1", "Windows Server", "m6id.xlarge",  "Shared Instances", "2", "40", "Hours/Week", "On-Demand", "General Purpose SSD (gp3)", "230", "2x Daily", "10")
)

for ($i = 0; $i -lt $controle.Length; $i++) {
    $row = $i + 2
    $data = $controle[$i]
    $ws2.Cells.Item($row, 1).Value = $data[0]                 # A Group
    $ws2.Cells.Item($row, 2).Value = $data[1]                 # B Description
    $ws2.Cells.Item($row, 3).Value = $data[2]                 # C AWS Region
    $ws2.Cells.Item($row, 4).Value = $data[3]                 # D Operating System
    $ws2.Cells.Item($row, 5).Value = $data[4]                 # E Instance Type
    $ws2.Cells.Item($row, 6).Value = $data[5]                 # F Tenancy
    Set-TextValue $ws2.Cells.Item($row, 7)  $data[6]          # G Number of Instances (text)
    Set-TextValue $ws2.Cells.Item($row, 8)  $data[7]          # H Assumed Usage (text)
    $ws2.Cells.Item($row, 9).Value = $data[8]                 # I Usage Type
    $ws2.Cells.Item($row, 10).Value = $data[9]                # J Purchasing Option
    $ws2.Cells.Item($row, 11).Value = $data[10]               # K Storage Type
    Set-TextValue $ws2.Cells.Item($row, 12) $data[11]         # L Storage amount (text)
    $ws2.Cells.Item($row, 15).Value = $data[12]                # O Snapshot Frequency
    Set-TextValue $ws2.Cells.Item($row, 16) $data[13]         # P EBS Snapshot amount (text)
}

# Drop the now-obsolete trailing row (previously Development / 2 x BACKUP 32GB...)
$ws2.Range("A6:P6").Delete()
